# Add a "get more" localization entry (time-trial "get more" button) as a
# new row in the Strings sheet, right after the "Challenge levels" row and
# before the "challenge title" row. Inserting the row pushes every
# following row down by one, which also pulls the rest of the
# localization table (and the implicitly-maintained shared-string table)
# into alignment with the target layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 89 (the old row 89 "challenge title" / "C H A L L E N G E"
# and everything below it shifts down to row 90+).
$ws.Rows("89:89").Insert()

# Populate the new row: column A holds the localization key, column B the
# English display text (same convention as every other row in the sheet).
$ws.Range("A89").Value = "get more"
$ws.Range("B89").Value = "Get more"

# Match the row height used by the other single-line rows in this block.
$ws.Rows("89:89").RowHeight = 13.4

# Reflect the new active cell (user selected the freshly-inserted row after
# editing), matching the saved view state.
[void]$ws.Range("A90").Select()
